$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 757, shifting existing rows 757+ down by one.
$ws.Rows.Item(757).Insert()

# Force the date-like text into the cell as plain text (not an auto-converted
# date serial number) by pre-formatting the cell as Text.
$ws.Cells.Item(757, 1).NumberFormat = "@"
$ws.Cells.Item(757, 2).NumberFormat = "@"

# Populate the newly inserted row 757 with the new data point.
$ws.Cells.Item(757, 1).Value = "2026/01/31"
$ws.Cells.Item(757, 2).Value = "土"
$ws.Cells.Item(757, 3).Value = 19
$ws.Cells.Item(757, 4).Value = 193

# Re-align the formatting of the new row with its neighbors (plain, unstyled
# cells) by copying the format from the row above over the Text override.
$ws.Cells.Item(756, 1).Copy()
$ws.Cells.Item(757, 1).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(756, 2).Copy()
$ws.Cells.Item(757, 2).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(756, 3).Copy()
$ws.Cells.Item(757, 3).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(756, 4).Copy()
$ws.Cells.Item(757, 4).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
